$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of values inserted at row 2; every prior row 2-10 shifts down into row 3-11.
# (Row 11's original data is discarded since there is no row 12 to push it into.)
$values = @(
    @(0.001217196552899444, 2.134347703298304, 20.89521853204074, 4.571128802827672, 4.673863460101999, 23),
    @(-0.08714184739790505, 2.431811455565427, 19.88066375958591, 4.458773795516645, 4.562828807598136, 22),
    @(-0.5723694534419483, 1.7269258304889, 9.411467503488769, 3.067811516943107, 3.088373965028863, 21),
    @(-0.1669402324491583, 1.364089054028298, 9.210978866053001, 3.034959450479199, 3.10908851507502, 20),
    @(-0.1633360453900471, 1.833324964313758, 12.49798808868182, 3.535249367255699, 3.628244747611019, 19),
    @(-0.2690867049885003, 1.819627553723605, 10.42808925027345, 3.229255216032553, 3.311319878293547, 18),
    @(-0.2014155535992504, 1.74506880832038, 10.82754791611375, 3.290523957687248, 3.385434377086962, 17),
    @(-0.09017090610257712, 1.776193825283892, 11.1038679904649, 3.332246688116727, 3.440269322504859, 16),
    @(-0.1106871006837081, 2.016227315994001, 13.2833031734348, 3.644626616463585, 3.77080678263623, 15),
    @(-0.1031043840192009, 1.961533181914308, 12.81524593535968, 3.579838814159051, 3.713433160618765, 14)
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = 2 + $i
    $rowVals = $values[$i]
    $ws.Cells.Item($row, 2).Value = $rowVals[0]
    $ws.Cells.Item($row, 3).Value = $rowVals[1]
    $ws.Cells.Item($row, 4).Value = $rowVals[2]
    $ws.Cells.Item($row, 5).Value = $rowVals[3]
    $ws.Cells.Item($row, 6).Value = $rowVals[4]
    $ws.Cells.Item($row, 7).Value = $rowVals[5]
}
